# LOM3218.xlsx - Build site at 2023-04-12 14:53:07 UTC
#
# Rewrites the course-plan worksheet: adds Portuguese "Objetivos", a second
# docente, Portuguese short/long syllabus text and a "Bibliografia" block,
# which pushes the remaining rows (Metodo/Criterio/Norma de
# recuperacao/Bibliografia) down and gives each of them its own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New text blocks introduced by this revision
# ---------------------------------------------------------------------
$objetivosPt = "Apresentar aos alunos ingressantes o entendimento do que seja a carreira e as bases conceituais da Engenharia Física."

$resumoPt = "A carreira de Engenharia Física. Conceitos básicos de Engenharia. Competências e habilidades de um engenheiro. Física conceitual. Realização de experimentos e projetos de Engenharia Física."

$programaPt = "A carreira de Engenharia Física. Cientistas x engenheiros: o papel interdisciplinar da Engenharia Física. Campos de atuação. A Física como ciência conceitual: Como aprender Física. Realização de demonstrações e experimentos científicos significativos de Física.Conceitos básicos de Engenharia. Habilidades e competências de um engenheiro.Desenvolvimento de um projeto temático de Engenharia Física.Competição entre projetos de diferentes grupos.Avaliação das competições e da disciplina como um todo."

$bibliografia = @"
ARAÚJO-MOREIRA, F. M. Engenharia Física: a Carreira do Novo Milênio, São Carlos: Gráfica e Editora Guillen & Andriolli, 2014.
BAZZO, A. B.; PEREIRA, L.T.V. Introdução à Engenharia. Editora da UFSC, Florianópolis, 1993.
ALEXANDER, C. K.; WATSON, J. A. Habilidades para uma carreira de sucesso na engenharia, Porto Alegre: AMGH Editora, 2015.
BROCKMAN, J. B. Introdução à Engenharia. LTC, Rio de Janeiro, 2009.
KNOWLEDGE FLOW. Engineering Physics - Ebook, Índia, 2015.
CHAVES, A. S.; VALADARES, E. C.; ALVES, E. G. Aplicações da Física Quântica do Transistor à Nanotecnologia, São Paulo: Livraria da Física, 2005.
"@

# ---------------------------------------------------------------------
# Row 10/11 ("Objetivos:" / "Objectives:") gains the new Portuguese text
# in column B/C (English text in row 11 is unchanged).
# ---------------------------------------------------------------------
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# ---------------------------------------------------------------------
# Row 12 ("Docentes responsáveis:") stays as-is. The two professors
# previously crammed into rows 13/14 alongside unrelated text now get
# their own clean rows with only B/C filled in.
# ---------------------------------------------------------------------
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(13).AutoFit()

$ws.Range("A14").ClearContents()
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Rows.Item(14).AutoFit()

# ---------------------------------------------------------------------
# Row 15 "Programa resumido:" now gets its own Portuguese summary
# (previously shared the English "Objetivos" text by mistake).
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = $resumoPt
$ws.Range("C15").Value = $resumoPt
$ws.Rows.Item(15).RowHeight = 60

# Row 16 "Short syllabus:" keeps its English text, only the row height shrinks.
$ws.Rows.Item(16).RowHeight = 60

# ---------------------------------------------------------------------
# Row 17 "Programa:" gains the full Portuguese syllabus text.
# ---------------------------------------------------------------------
$ws.Range("B17").Value = $programaPt
$ws.Range("C17").Value = $programaPt
$ws.Rows.Item(17).RowHeight = 120

# Row 18 "Syllabus:" (English) shifts down from the old row 17 slot.
$ws.Range("A18").Value = "Syllabus:"
$ws.Range("B18").Value = "The career of Engineering Physics. Scientists x engineers: the interdisciplinary role of Engineering Physics. Fields of action.Physics as a conceptual science: How to learn Physics. Realization of demonstrations and significant scientific experiments in Physics.Basic engineering concepts. Skills and competences of an engineer.Development of a thematic project of Physical Engineering.Competition between projects from different groups.Evaluation of competitions and the discipline as a whole."
$ws.Range("C18").Value = "The career of Engineering Physics. Scientists x engineers: the interdisciplinary role of Engineering Physics. Fields of action.Physics as a conceptual science: How to learn Physics. Realization of demonstrations and significant scientific experiments in Physics.Basic engineering concepts. Skills and competences of an engineer.Development of a thematic project of Physical Engineering.Competition between projects from different groups.Evaluation of competitions and the discipline as a whole."
$ws.Rows.Item(18).RowHeight = 120

# ---------------------------------------------------------------------
# Row 19 becomes the standalone "Avaliação:" header (no B/C content).
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Avaliação:"
$ws.Range("B19:C19").ClearContents()
$ws.Rows.Item(19).AutoFit()

# ---------------------------------------------------------------------
# Row 20 "Método:" keeps its own text (already correct) - just confirm.
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "As atividades práticas e os projetos que serão desenvolvidos durante as aulas serão avaliados por docentes e pelos alunos (processo de avaliação crítica)."
$ws.Range("C20").Value = "As atividades práticas e os projetos que serão desenvolvidos durante as aulas serão avaliados por docentes e pelos alunos (processo de avaliação crítica)."
$ws.Rows.Item(20).RowHeight = 60

# Row 21 "Critério:" moves up from the old "Norma de recuperação" text slot.
$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "A média final será uma composição de fatores relativos à participação do aluno nos trabalhos desenvolvidos, conjuntamente com o rendimento de seu grupo."
$ws.Range("C21").Value = "A média final será uma composição de fatores relativos à participação do aluno nos trabalhos desenvolvidos, conjuntamente com o rendimento de seu grupo."
$ws.Rows.Item(21).RowHeight = 60

# New row 22: "Norma de recuperação:"
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "Devido às características da disciplina, não será oferecida recuperação."
$ws.Range("C22").Value = "Devido às características da disciplina, não será oferecida recuperação."
$ws.Rows.Item(22).RowHeight = 60

# New row 23: "Bibliografia:" with the new reference list.
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = $bibliografia
$ws.Range("C23").Value = $bibliografia
$ws.Rows.Item(23).RowHeight = 120
